$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 850.1429000000001
$ws.Range("I99").Value = 575.2
$ws.Range("J99").Value = 1537.5
$ws.Range("K99").Value = 1725.6
$ws.Range("L99").Value = 4612.5
$ws.Range("M99").Value = -227.6000000000001
$ws.Range("N99").Value = -7608.5
$ws.Range("H112").Value = 1058.4286
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 1063.9062
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 3191.7186
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -5407.7186
$ws.Range("H116").Value = 7031.136
$ws.Range("I116").Value = 10123.75
$ws.Range("J116").Value = 3320
$ws.Range("K116").Value = 10123.75
$ws.Range("L116").Value = 3320
$ws.Range("M116").Value = -6681.75
$ws.Range("N116").Value = -10204
$ws.Range("H129").Value = 956.9383
$ws.Range("I129").Value = 746.1429000000001
$ws.Range("J129").Value = 976.87836
$ws.Range("K129").Value = 2238.4287
$ws.Range("L129").Value = 2930.63508
$ws.Range("M129").Value = 2761.5713
$ws.Range("N129").Value = -12930.63508
$ws.Range("H132").Value = 1100.8276
$ws.Range("I132").Value = 849.06384
$ws.Range("K132").Value = 2547.19152
$ws.Range("M132").Value = -17.19152000000031
$ws.Range("H138").Value = 2505.486
$ws.Range("I138").Value = 1154.8158
$ws.Range("J138").Value = 4015.0588
$ws.Range("K138").Value = 3464.4474
$ws.Range("L138").Value = 12045.1764
$ws.Range("M138").Value = 1675.5526
$ws.Range("N138").Value = -22325.1764

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3327.3264
$ws.Range("I61").Value = 3746.8684
$ws.Range("J61").Value = 1878
$ws.Range("K61").Value = 3746.8684
$ws.Range("L61").Value = 1878
$ws.Range("M61").Value = -3534.8684
$ws.Range("N61").Value = -2302
$ws.Range("H74").Value = 1512.7307
$ws.Range("I74").Value = 1410.1666
$ws.Range("J74").Value = 1743.5
$ws.Range("K74").Value = 1410.1666
$ws.Range("L74").Value = 1743.5
$ws.Range("M74").Value = -536.1666
$ws.Range("N74").Value = -3491.5
$ws.Range("H77").Value = 1512.7307
$ws.Range("I77").Value = 1410.1666
$ws.Range("J77").Value = 1743.5
$ws.Range("K77").Value = 7050.833000000001
$ws.Range("L77").Value = 8717.5
$ws.Range("M77").Value = -2682.833000000001
$ws.Range("N77").Value = -17453.5
$ws.Range("H132").Value = 1668969.5
$ws.Range("I132").Value = 1807.6316
$ws.Range("J132").Value = 4548613
$ws.Range("K132").Value = 5422.8948
$ws.Range("L132").Value = 13645839
$ws.Range("M132").Value = -2892.8948
$ws.Range("N132").Value = -13650899
$ws.Range("H136").Value = 3327.3264
$ws.Range("I136").Value = 3746.8684
$ws.Range("J136").Value = 1878
$ws.Range("K136").Value = 11240.6052
$ws.Range("L136").Value = 5634
$ws.Range("M136").Value = -8690.6052
$ws.Range("N136").Value = -10734

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3454.32
$ws.Range("I134").Value = 3732.925
$ws.Range("J134").Value = 2339.9
$ws.Range("K134").Value = 11198.775
$ws.Range("L134").Value = 7019.700000000001
$ws.Range("M134").Value = -8663.775000000001
$ws.Range("N134").Value = -12089.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15255.739
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 15255.739
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H58").Value = 1027.695
$ws.Range("I58").Value = 598.04443
$ws.Range("J58").Value = 2408.7144
$ws.Range("K58").Value = 598.04443
$ws.Range("L58").Value = 2408.7144
$ws.Range("M58").Value = -395.04443
$ws.Range("N58").Value = -2814.7144
$ws.Range("H132").Value = 1383.4445
$ws.Range("I132").Value = 1020.04254
$ws.Range("J132").Value = 2450.9375
$ws.Range("K132").Value = 3060.12762
$ws.Range("L132").Value = 7352.8125
$ws.Range("M132").Value = -530.1276200000002
$ws.Range("N132").Value = -12412.8125
$ws.Range("H134").Value = 1942.3115
$ws.Range("I134").Value = 2339.1943
$ws.Range("J134").Value = 1370.8
$ws.Range("K134").Value = 7017.5829
$ws.Range("L134").Value = 4112.4
$ws.Range("M134").Value = -4482.5829
$ws.Range("N134").Value = -9182.4
$ws.Range("H136").Value = 1027.695
$ws.Range("I136").Value = 598.04443
$ws.Range("J136").Value = 2408.7144
$ws.Range("K136").Value = 1794.13329
$ws.Range("L136").Value = 7226.1432
$ws.Range("M136").Value = 755.8667099999998
$ws.Range("N136").Value = -12326.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 152.375
$ws.Range("J11").Value = 180
$ws.Range("L11").Value = 540
$ws.Range("N11").Value = -820
$ws.Range("H75").Value = 20409748
$ws.Range("J75").Value = 20409748
$ws.Range("L75").Value = 61229244
$ws.Range("N75").Value = -61231240
$ws.Range("H78").Value = 20409748
$ws.Range("J78").Value = 20409748
$ws.Range("L78").Value = 183687732
$ws.Range("N78").Value = -183697716
$ws.Range("H117").Value = 22233230
$ws.Range("I117").Value = 14632.571
$ws.Range("J117").Value = 41674504
$ws.Range("K117").Value = 43897.713
$ws.Range("L117").Value = 125023512
$ws.Range("M117").Value = -40455.713
$ws.Range("N117").Value = -125030396
$ws.Range("H121").Value = 1035.48
$ws.Range("I121").Value = 590
$ws.Range("J121").Value = 1120.3334
$ws.Range("K121").Value = 1770
$ws.Range("L121").Value = 3361.0002
$ws.Range("M121").Value = -460
$ws.Range("N121").Value = -5981.0002
$ws.Range("H129").Value = 17544946
$ws.Range("J129").Value = 1575.5555
$ws.Range("L129").Value = 4726.666499999999
$ws.Range("N129").Value = -14726.6665
$ws.Range("H130").Value = 3827.5
$ws.Range("I130").Value = 515
$ws.Range("J130").Value = 7140
$ws.Range("K130").Value = 1545
$ws.Range("L130").Value = 21420
$ws.Range("M130").Value = 3475
$ws.Range("N130").Value = -31460
$ws.Range("H131").Value = 1563403.4
$ws.Range("I131").Value = 5882996.5
$ws.Range("K131").Value = 17648989.5
$ws.Range("M131").Value = -17643949.5
$ws.Range("H136").Value = 17868074
$ws.Range("J136").Value = 4416.6665
$ws.Range("L136").Value = 13249.9995
$ws.Range("N136").Value = -23449.9995
$ws.Range("H139").Value = 10207947
$ws.Range("I139").Value = 20004812
$ws.Range("J139").Value = 2880.6667
$ws.Range("K139").Value = 60014436
$ws.Range("L139").Value = 8642.000100000001
$ws.Range("M139").Value = -60009296
$ws.Range("N139").Value = -18922.0001
$ws.Range("H140").Value = 20834788
$ws.Range("I140").Value = 21740584
$ws.Range("J140").Value = 1500
$ws.Range("K140").Value = 65221752
$ws.Range("L140").Value = 4500
$ws.Range("M140").Value = -65216572
$ws.Range("N140").Value = -14860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 18025
$ws.Range("J57").Value = 18025
$ws.Range("L57").Value = 18025
$ws.Range("N57").Value = -19665
$ws.Range("H132").Value = 1709.878
$ws.Range("I132").Value = 1317.5714
$ws.Range("J132").Value = 3998.3333
$ws.Range("K132").Value = 3952.7142
$ws.Range("L132").Value = 11994.9999
$ws.Range("M132").Value = -1422.7142
$ws.Range("N132").Value = -17054.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 20000700
$ws.Range("I93").Value = 619.125
$ws.Range("J93").Value = 55556396
$ws.Range("K93").Value = 619.125
$ws.Range("L93").Value = 55556396
$ws.Range("M93").Value = 628.875
$ws.Range("N93").Value = -55558892
$ws.Range("H132").Value = 13363029
$ws.Range("I132").Value = 20557238
$ws.Range("J132").Value = 2355.0715
$ws.Range("K132").Value = 61671714
$ws.Range("L132").Value = 7065.2145
$ws.Range("M132").Value = -61669184
$ws.Range("N132").Value = -12125.2145
$ws.Range("H136").Value = 4481.7896
$ws.Range("I136").Value = 3868.1555
$ws.Range("K136").Value = 11604.4665
$ws.Range("M136").Value = -9054.466499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18076.482
$ws.Range("I132").Value = 25042.512
$ws.Range("J132").Value = 1276.0588
$ws.Range("K132").Value = 75127.53599999999
$ws.Range("L132").Value = 3828.1764
$ws.Range("M132").Value = -72597.53599999999
$ws.Range("N132").Value = -8888.1764
$ws.Range("H136").Value = 6668700
$ws.Range("I136").Value = 2194.32
$ws.Range("J136").Value = 20001710
$ws.Range("K136").Value = 6582.960000000001
$ws.Range("L136").Value = 60005130
$ws.Range("M136").Value = -4032.960000000001
$ws.Range("N136").Value = -60010230
